$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two bold runs "Requerimientos " + "Reconocidos" into a single
#    run "Requerimientos Reconocidos" (title paragraph at top of document).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Requerimientos Reconocidos", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Requerimientos Reconocidos", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Table-wide geometry tweaks, applied to both tables in the document:
#    - tblInd: 170 -> 160 dxa (8.5pt -> 8.0pt)
#    - left cell margin (tcMar / per cell): 90 -> 80 dxa (4.5pt -> 4.0pt)
# ---------------------------------------------------------------------------
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)

    # table indent from margin
    $tbl.Rows.LeftIndent = 8.0

    # table-wide default left cell margin (affects w:tblCellMar)
    $tbl.LeftPadding = 4.0

    # explicit per-cell left margin overrides (w:tcMar) on every existing cell
    foreach ($row in $tbl.Rows) {
        foreach ($cell in $row.Cells) {
            $cell.LeftPadding = 4.0
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Column width tweaks on table 1 only:
#    - column 2: 3644 -> 3643 dxa (182.2pt -> 182.15pt)
#    - column 3: 4621 -> 4622 dxa (231.05pt -> 231.1pt)
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Columns.Item(2).Width = 182.15
$t1.Columns.Item(3).Width = 231.1

# ---------------------------------------------------------------------------
# 4) Add a new paragraph "CAMBIO 1" under "Administrar productos" in the
#    requirement table's row for R1 (row 2, column 2).
# ---------------------------------------------------------------------------
$t1b = $d.Tables.Item(1)
$reqCell = $t1b.Cell(2, 2)
$reqCell.Range.InsertAfter([char]13 + "CAMBIO 1")

# ---------------------------------------------------------------------------
# 5) Normal style's language tag re-normalizes to w:val, w:eastAsia, w:bidi
#    order when the style is touched; re-assert the (unchanged) language to
#    trigger that canonical re-serialization.
# ---------------------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.Font.LanguageID = "es-EC"
